$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# 1) Introduction paragraph: drop the leading empty paragraph, indent the
#    text paragraph and rewrite its content (split in two runs around a
#    relocated "_GoBack" bookmark).
# ----------------------------------------------------------------------

$pEmpty = $d.Paragraphs(6)
$pEmpty.Range.Delete()

$pIntro = $d.Paragraphs(6)
$pIntro.LeftIndent = 36

$introRange = $pIntro.Range
$introRange.MoveEnd(1, -1)
$introStart = $introRange.Start

$part1 = "Este documento contém o planejamento geral do projeto do Sistema Industrial que será desenvolvido para auxiliar na gestão das empresas que trabalham com produção de produtos, com a finalidade de um total "
$part2 = "controle e diminuição de gastos extras. "
$introRange.Text = $part1 + $part2

# Move the "_GoBack" bookmark from its old location (end of the document
# body) to the point between the two halves of the text we just typed.
$d.Bookmarks("_GoBack").Delete()
$bmPoint = $d.Range($introStart + $part1.Length, $introStart + $part1.Length)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# ----------------------------------------------------------------------
# 2) Requirements table: merge runs that were only split apart because of
#    a (now stale) lastRenderedPageBreak, and drop the remaining stale
#    lastRenderedPageBreak markers.
# ----------------------------------------------------------------------

$d.Content.Find.Execute("Incluir componentes no produto final", $true, $false, $false, $false, $false, $true, 1, $false, "Incluir componentes no produto final", 2) | Out-Null
$d.Content.Find.Execute("Inclusão de matéria prima e quantidade necessária para a produção de cada produto final", $true, $false, $false, $false, $false, $true, 1, $false, "Inclusão de matéria prima e quantidade necessária para a produção de cada produto final", 2) | Out-Null

Write-Output "done"
